$d = $word.ActiveDocument

# Номер аудитории: 42 -> 43
$d.Content.Find.Execute("42", $true, $false, $false, $false, $false,
                         $true, 1, $false, "43", 2)

# Отчетный период: с 18.05.2024 по 01.06.2024 -> с 03.06.2024 по 13.06.2024
$d.Content.Find.Execute("с 18.05.2024 по 01.06.2024", $true, $false, $false, $false, $false,
                         $true, 1, $false, "с 03.06.2024 по 13.06.2024", 2)

# Table cell: dfgs -> ПК №1
$d.Content.Find.Execute("dfgs", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ПК №1", 2)

# Table cell: sdfgdsfg -> Установка новые драйвера
$d.Content.Find.Execute("sdfgdsfg", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Установка новые драйвера", 2)
